$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.491.98"
$ws.Cells.Item(2, 5).Value = "  -0.39%  "
$ws.Cells.Item(3, 4).Value = "1.820.46"
$ws.Cells.Item(3, 5).Value = "  -0.33%  "
$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$ws.Cells.Item(5, 4).Value = "'317.32"
$ws.Cells.Item(5, 5).Value = "  +0.20%  "
$ws.Cells.Item(7, 4).Value = "'0.5166"
$ws.Cells.Item(7, 5).Value = "  -2.72%  "
$ws.Cells.Item(8, 5).Value = "  -2.09%  "
$ws.Cells.Item(9, 4).Value = "'0.08462"
$ws.Cells.Item(9, 5).Value = "  +9.04%  "
$ws.Cells.Item(10, 4).Value = "'41.82"
$ws.Cells.Item(10, 5).Value = "  -0.56%  "
$ws.Cells.Item(11, 5).Value = "  -0.63%  "
$ws.Cells.Item(12, 4).Value = "'6.446"
$ws.Cells.Item(12, 5).Value = "  +1.98%  "
$ws.Cells.Item(13, 4).Value = "'21.05"
$ws.Cells.Item(13, 5).Value = "  -0.35%  "
$ws.Cells.Item(14, 4).Value = "'1.003"
$ws.Cells.Item(14, 5).Value = "  +0.02%  "
$ws.Cells.Item(15, 4).Value = "'7.517"
$ws.Cells.Item(15, 5).Value = "  -0.66%  "
$ws.Cells.Item(16, 4).Value = "1.814.53"
$ws.Cells.Item(16, 5).Value = "  -2.53%  "
$ws.Cells.Item(17, 4).Value = "'0.00001143"
$ws.Cells.Item(17, 5).Value = "  +4.95%  "
$ws.Cells.Item(18, 4).Value = "'92.84"
$ws.Cells.Item(18, 5).Value = "  -0.46%  "
$ws.Cells.Item(19, 4).Value = "'0.06633"
$ws.Cells.Item(19, 5).Value = "  +0.15%  "
$ws.Cells.Item(20, 4).Value = "'17.76"
$ws.Cells.Item(20, 5).Value = "  -0.12%  "
$ws.Cells.Item(22, 4).Value = "'6.092"
$ws.Cells.Item(22, 5).Value = "  +0.14%  "
$ws.Cells.Item(23, 4).Value = "28.525.21"
$ws.Cells.Item(23, 5).Value = "  -0.30%  "
$ws.Cells.Item(24, 4).Value = "'11.47"
$ws.Cells.Item(24, 5).Value = "  +2.40%  "
$ws.Cells.Item(25, 5).Value = "  +1.56%  "
$ws.Cells.Item(26, 4).Value = "'21.07"
$ws.Cells.Item(26, 5).Value = "  +1.29%  "
$ws.Cells.Item(27, 4).Value = "'159.28"
$ws.Cells.Item(27, 5).Value = "  +1.50%  "
$ws.Cells.Item(28, 4).Value = "2.025.59"
$ws.Cells.Item(28, 5).Value = "  +0.06%  "
$ws.Cells.Item(29, 4).Value = "'2.396"
$ws.Cells.Item(29, 5).Value = "  -0.88%  "
$ws.Cells.Item(30, 4).Value = "'125.53"
$ws.Cells.Item(30, 5).Value = "  +0.20%  "
$ws.Cells.Item(31, 4).Value = "'0.1090"
$ws.Cells.Item(31, 5).Value = "  -3.46%  "
$ws.Cells.Item(32, 5).Value = "  -4.80%  "
$ws.Cells.Item(33, 4).Value = "'5.728"
$ws.Cells.Item(33, 5).Value = "  -0.32%  "
$ws.Cells.Item(34, 4).Value = "'0.07444"
$ws.Cells.Item(34, 5).Value = "  +1.60%  "
$ws.Cells.Item(35, 4).Value = "'3.647"
$ws.Cells.Item(35, 5).Value = "  -0.45%  "
$ws.Cells.Item(36, 4).Value = "'0.2233"
$ws.Cells.Item(36, 5).Value = "  -1.66%  "
$ws.Cells.Item(37, 5).Value = "  +0.44%  "
$ws.Cells.Item(38, 4).Value = "'5.210"
$ws.Cells.Item(38, 5).Value = "  +0.04%  "
$ws.Cells.Item(40, 4).Value = "'0.6314"
$ws.Cells.Item(40, 5).Value = "  +0.28%  "
$ws.Cells.Item(41, 4).Value = "'11.27"
$ws.Cells.Item(41, 5).Value = "  -1.32%  "
$ws.Cells.Item(42, 4).Value = "'1.190"
$ws.Cells.Item(42, 5).Value = "  -0.38%  "
$ws.Cells.Item(43, 5).Value = "  -0.04%  "
$ws.Cells.Item(44, 4).Value = "'13.59"
$ws.Cells.Item(44, 5).Value = "  +0.36%  "
$ws.Cells.Item(45, 4).Value = "'3.785"
$ws.Cells.Item(45, 5).Value = "  +1.74%  "
$ws.Cells.Item(46, 4).Value = "'0.5948"
$ws.Cells.Item(46, 5).Value = "  +0.04%  "
$ws.Cells.Item(47, 4).Value = "'126.45"
$ws.Cells.Item(47, 5).Value = "  +0.66%  "
$ws.Cells.Item(48, 4).Value = "'1.988"
$ws.Cells.Item(48, 5).Value = "  -0.60%  "
$ws.Cells.Item(49, 5).Value = "  +0.77%  "
$ws.Cells.Item(50, 4).Value = "'0.06983"
$ws.Cells.Item(50, 5).Value = "  +0.29%  "
$ws.Cells.Item(51, 4).Value = "'74.37"
$ws.Cells.Item(51, 5).Value = "  -0.43%  "
